# Applies the "user mode completed and business mode start" update:
#  - Social sheet row1 (A1/B1) repurposed from the "pragas" e-mail hyperlink
#    row into a phone-number / password row (A1 becomes a plain number with
#    left/top aligned formatting, B3 becomes the Synctag password), the
#    other rows shuffle up to the next account, and a few sheets pick up a
#    refreshed cell selection.

$wb = $excel.ActiveWorkbook

$wsLog    = $wb.Worksheets.Item("SynctagLog")
$wsOutlook = $wb.Worksheets.Item("OutLook")
$wsSocial = $wb.Worksheets.Item("Social")
$wsMedia  = $wb.Worksheets.Item("Media")

# --- Social sheet: row values shift / row1 repurposed -----------------
$wsSocial.Range("A1").Value = 7639416734
$wsSocial.Range("A1").HorizontalAlignment = -4131
$wsSocial.Range("A1").VerticalAlignment = -4160

$wsSocial.Range("A2").Value = "ranjithkumar.hashinn@gmail.com"
$wsSocial.Range("A3").Value = "prathapkumar.hashinn@gmail.com"
$wsSocial.Range("B3").Value = "Synctag@1"

# A1's hyperlink keeps pointing at pragas.hashinn@gmail.com, but the cell
# text is now the phone number, so preserve the old label as the display
# text. A3's hyperlink now matches its cell text exactly, so its stale
# display override is cleared.
foreach ($hl in $wsSocial.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$A`$1") {
        $hl.TextToDisplay = "pragas.hashinn@gmail.com"
    } elseif ($addr -eq "`$A`$3") {
        $hl.TextToDisplay = ""
    }
}

# --- Selection bookmarks ------------------------------------------------
$wsLog.Range("A2").Select()
$wsSocial.Range("C17").Select()

# Re-activate SynctagLog so it remains the selected tab, matching its
# original tabSelected sheetView state.
$wsLog.Range("A2").Select()
